$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 294, shifting existing rows 294..358 down to 295..359.
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new price observation.
$ws.Cells.Item(294, 1).Value  = 3
$ws.Cells.Item(294, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(294, 3).Value  = "Coquimbo"
$ws.Cells.Item(294, 4).Value  = 44511
$ws.Cells.Item(294, 5).Value  = 5
$ws.Cells.Item(294, 6).Value  = 100112045
$ws.Cells.Item(294, 7).Value  = "Zapallo"
$ws.Cells.Item(294, 8).Value  = "Camote"
$ws.Cells.Item(294, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(294, 10).Value = 160
$ws.Cells.Item(294, 11).Value = 600
$ws.Cells.Item(294, 12).Value = 600
$ws.Cells.Item(294, 13).Value = 600
$ws.Cells.Item(294, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(294, 15).Value = "Perú"
$ws.Cells.Item(294, 16).Value = 600
$ws.Cells.Item(294, 17).Value = 1
$ws.Cells.Item(294, 18).Value = "Hortaliza"
